# Bytter rekkefølge på matte-operatorer
#
# This script edits the "Matematikk-operatorer" table on slide 1 so the
# Divisjon/Multiplikasjon rows swap places, and the Eksponent/Heltalldivisjon
# rows swap places (with a couple of small label tweaks), and updates two
# labels in the "Matematiske funksjoner" table.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the two tables we need by inspecting their header cell text,
# since shape ordering should be stable but we verify defensively.
$opTable = $null
$funcTable = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $hdr = $sh.Table.Cell(1, 1).Shape.TextFrame.TextRange.Text
        if ($hdr -eq "Matematikk-operatorer") {
            $opTable = $sh.Table
        } elseif ($hdr -eq "Matematiske funksjoner") {
            $funcTable = $sh.Table
        }
    }
}

# --- "Matematikk-operatorer" table -----------------------------------
# Row layout (1-indexed): 1=header, 2=+, 3=-, 4=/, 5=*, 6=**, 7=%, 8=//

# Row 4 (was "/" Divisjon) now becomes "*" Multiplikasjon
$opTable.Cell(4, 1).Shape.TextFrame.TextRange.Text = "*"
$opTable.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Multiplikasjon`rEksempel:`r9 * 4 = 36"

# Row 5 (was "*" Multiplikasjon) now becomes "/" Divisjon
$opTable.Cell(5, 1).Shape.TextFrame.TextRange.Text = "/"
$opTable.Cell(5, 2).Shape.TextFrame.TextRange.Text = "Divisjon`rEksempel:`r9 / 4 = 2.25"

# Row 6 (was "**" Eksponent) now becomes "//" Heltalldivisjon
$opTable.Cell(6, 1).Shape.TextFrame.TextRange.Text = "//"
$opTable.Cell(6, 2).Shape.TextFrame.TextRange.Text = "Heltalldivisjon`rEksempel:`r9 // 4 = 2"

# Row 7 ("%" Modulus / rest) keeps its operator/example, only the label
# "Eksempel" gains a trailing colon.
$opTable.Cell(7, 2).Shape.TextFrame.TextRange.Paragraphs(2, 1).Text = "Eksempel:"

# Row 8 (was "//" Heltalldivisjon) now becomes "**" Eksponent, including
# the trailing blank paragraph that the original Eksponent row had.
$opTable.Cell(8, 1).Shape.TextFrame.TextRange.Text = "**"
$opTable.Cell(8, 2).Shape.TextFrame.TextRange.Text = "Eksponent`rEksempel`r9 ** 4 = 6561`r"

# --- "Matematiske funksjoner" table -----------------------------------
# Row 7: math.ceil(x) -> "Runder av opp til heltall"
$funcTable.Cell(7, 2).Shape.TextFrame.TextRange.Text = "Runder av opp til heltall"
# Row 8: math.floor(x) -> "Runder av ned til heltall"
$funcTable.Cell(8, 2).Shape.TextFrame.TextRange.Text = "Runder av ned til heltall"
